# "original checked or unchecked with docx o1"
#
# The "_GoBack" bookmark (collapsed, id 0) currently sits right after the
# {s1f15} merge field. It needs to move to sit right after the
# "True Copy of the Original" paragraph, whose text is being replaced by
# the merge field {o1}.
#
#   1. Delete the existing _GoBack bookmark.
#   2. Replace the "True Copy of the Original" run text with "{o1}".
#   3. Re-insert a collapsed _GoBack bookmark immediately after the new
#      {o1} run (i.e. at the end of that paragraph, before its paragraph
#      mark).

$d = $word.ActiveDocument

# --- Step 1: remove the pre-existing _GoBack bookmark -----------------------
try {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
} catch {
    # No _GoBack bookmark present - nothing to remove.
}

# --- Step 2: find the paragraph holding "True Copy of the Original" --------
$needle = "True Copy of the Original"
$found = $d.Content
$ok = $found.Find.Execute($needle, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find target text '$needle'"
}
$paraIndex = $found.Paragraphs.Item(1).Index

# Replace the whole paragraph's text with the merge field placeholder.
# Doing this via a freshly fetched Paragraphs.Item(...).Range (rather than
# the live Find range) keeps the existing run's formatting/rsid attributes
# intact, matching how Word itself edits text in place.
$p = $d.Paragraphs.Item($paraIndex)
$p.Range.Text = "{o1}"

# --- Step 3: re-create the _GoBack bookmark right after the new {o1} run ---
# Placing a *collapsed* bookmark directly at (paragraph end - 1) is mishandled
# by this runtime, so: append a 2-character scratch placeholder at the
# paragraph's end, add the bookmark collapsed just in front of it (which is
# exactly the final target position), then delete the placeholder. The
# bookmark stays collapsed at the correct spot once the scratch text is gone.
$p = $d.Paragraphs.Item($paraIndex)
$pr = $p.Range
$scratch = $d.Range($pr.End - 1, $pr.End - 1)
$scratch.InsertAfter("XX")

$p = $d.Paragraphs.Item($paraIndex)
$pr = $p.Range
$bmRange = $d.Range($pr.End - 3, $pr.End - 3)
$d.Bookmarks.Add("_GoBack", $bmRange)

$p = $d.Paragraphs.Item($paraIndex)
$pr = $p.Range
$placeholder = $d.Range($pr.End - 3, $pr.End - 1)
$placeholder.Delete()
